$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cells to reflect 2020-12-16 data refresh for "Fonds de solidarite" volet 2
# Each cell is forced to Text format ("@") before assigning the value so that
# the numeric-looking strings are preserved exactly (matching the original inlineStr cells)
# rather than being converted to floating point numbers.

$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "496"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "441"

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2890139.57"

$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "227"

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "2008003.64"

$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "70"

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "1021337.89"

$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "794"

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "3815114.93"

$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "385"

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "3324428.91"

$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "150"

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "2102100.47"

$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "61"

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "1221474.34"

$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "8"

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "189000.00"

$ws.Range("C56").NumberFormat = "@"
$ws.Range("C56").Value = "1026"

$ws.Range("E56").NumberFormat = "@"
$ws.Range("E56").Value = "5835085.93"

$ws.Range("C57").NumberFormat = "@"
$ws.Range("C57").Value = "512"

$ws.Range("E57").NumberFormat = "@"
$ws.Range("E57").Value = "4821219.81"

$ws.Range("C58").NumberFormat = "@"
$ws.Range("C58").Value = "189"

$ws.Range("E58").NumberFormat = "@"
$ws.Range("E58").Value = "2070594.11"

$ws.Range("C59").NumberFormat = "@"
$ws.Range("C59").Value = "68"

$ws.Range("E59").NumberFormat = "@"
$ws.Range("E59").Value = "1055181.06"

$ws.Range("C63").NumberFormat = "@"
$ws.Range("C63").Value = "5714"

$ws.Range("E63").NumberFormat = "@"
$ws.Range("E63").Value = "24444042.70"

$ws.Range("C64").NumberFormat = "@"
$ws.Range("C64").Value = "3147"

$ws.Range("E64").NumberFormat = "@"
$ws.Range("E64").Value = "19420486.38"

$ws.Range("C67").NumberFormat = "@"
$ws.Range("C67").Value = "50"

$ws.Range("E67").NumberFormat = "@"
$ws.Range("E67").Value = "1842881.18"

$ws.Range("C74").NumberFormat = "@"
$ws.Range("C74").Value = "4"

$ws.Range("D74").NumberFormat = "@"
$ws.Range("D74").Value = "4"

$ws.Range("E74").NumberFormat = "@"
$ws.Range("E74").Value = "140000.00"
